# Commit: "identifiers for org, pract, patient..."
#
# Refresh the publish metadata on the "Metadata" sheet of the ValueSet
# workbook and record the resource's jurisdiction:
#   - Date            -> new publish timestamp
#   - Publisher       -> credit "D Foulkes"
#   - Contact         -> credit "D Foulkes"
#   - (new) Jurisdiction / Australia, inserted right after Contact

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Update Date / Publisher / Contact values (rows 8-10). ---------------
$ws.Range("B8").Value = "2024-05-20T17:01:27+10:00"
$ws.Range("B9").Value = "D Foulkes - Northern Australia Regional Digital Health Collaborative"
$ws.Range("B10").Value = "D Foulkes - Northern Australia Regional Digital Health Collaborative (https://nardhc.org)"

# --- Make room for a new "Jurisdiction" row right after "Contact". -------
# The sheet currently ends at row 14 (A1:B14); grow it to row 15 first by
# copying the formatting of the last existing row down into the new one,
# so the extra row matches the look of the rest of the table.
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Snapshot the current "Description" .. "Immutable" block (rows 11-14)
# before shifting it down one row to make room for "Jurisdiction".
$a11 = $ws.Range("A11").Value()
$b11 = $ws.Range("B11").Value()
$a12 = $ws.Range("A12").Value()
$b12 = $ws.Range("B12").Value()
$a13 = $ws.Range("A13").Value()
$b13 = $ws.Range("B13").Value()
$a14 = $ws.Range("A14").Value()
$b14 = $ws.Range("B14").Value()

$ws.Range("A15").Value = $a14
$ws.Range("B15").Value = $b14
$ws.Range("A14").Value = $a13
$ws.Range("B14").Value = $b13
$ws.Range("A13").Value = $a12
$ws.Range("B13").Value = $b12
$ws.Range("A12").Value = $a11
$ws.Range("B12").Value = $b11

# --- Write the new Jurisdiction / Australia row. --------------------------
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "Australia"
